$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B5").Formula = "=TEXT(1,""0"")"
$ws.Range("B5").Copy()
$ws.Range("B11").PasteSpecial(-4163)
$ws.Range("B5").ClearContents()
